$wb = $excel.ActiveWorkbook

# Rename worksheets (tab names) to the new task-order identifiers
$wb.Worksheets.Item(1).Name = "GNG_TO-1650291252503433"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912546333199"
$wb.Worksheets.Item(3).Name = "RS_TO-16502912546355336"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912546980171"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502912547939196"

# Sheet 1 (GNG) - update stim file names
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912524744043.csv"
$ws1.Range("B3").Value = "GNG_stims-16502912524874024.csv"
$ws1.Range("B4").Value = "go_stims-16502912524884043.csv"
$ws1.Range("B5").Value = "GNG_stims-1650291252501436.csv"

# Sheet 2 (NB) - update stim file names
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "OB-16502912537769911.csv"
$ws2.Range("B3").Value = "ZB-match_3-1650291252689762.csv"
$ws2.Range("B4").Value = "OB-16502912528607442.csv"
$ws2.Range("B5").Value = "ZB-match_1-16502912527437484.csv"
$ws2.Range("B6").Value = "TB-16502912546106806.csv"
$ws2.Range("B7").Value = "TB-16502912540160923.csv"
$ws2.Range("B8").Value = "OB-16502912531317453.csv"
$ws2.Range("B9").Value = "ZB-match_5-16502912525364084.csv"
$ws2.Range("B10").Value = "TB-16502912541499639.csv"

# Sheet 3 (RS) - no content changes

# Sheet 4 (TOL) - update stim file names
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912546492429.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912546386163.csv"
$ws4.Range("B4").Value = "MM_stims-16502912546810765.csv"
$ws4.Range("B5").Value = "ZM_stims-16502912546502457.csv"
$ws4.Range("B6").Value = "MM_stims-16502912546970005.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912546820772.csv"

# Sheet 5 (vSAT) - update stim file names
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-1650291254768874.csv"
$ws5.Range("B3").Value = "SAT_stims-16502912547271347.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502912547428563.csv"
$ws5.Range("B5").Value = "SAT_stims-16502912547029545.csv"
